# Auto-generated edit script.
# Applies row-content updates to rows 13-33 of sheet "Artfynd" in before.xlsx
# (species-occurrence records were reshuffled among rows; two rows gained new
# occurrence ids with slightly different observation times).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("A13").Value = 111574338
$ws.Range("B13").Value = 89686
$ws.Range("E13").Value = 658
$ws.Range("F13").Value = "Rosenticka"
$ws.Range("G13").Value = "Rhodofomes roseus"
$ws.Range("H13").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("I13").Value = ""
$ws.Range("Q13").Value = 562557.3535548041
$ws.Range("R13").Value = 6954757.635990249
$ws.Range("AC13").Value = ""

# Row 14
$ws.Range("A14").Value = 111574334
$ws.Range("B14").Value = 89405
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 1202
$ws.Range("F14").Value = "Ullticka"
$ws.Range("G14").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H14").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q14").Value = 562557.3535548041
$ws.Range("R14").Value = 6954757.635990249
$ws.Range("Z14").Value = "15:26"
$ws.Range("AB14").Value = "15:26"

# Row 15
$ws.Range("A15").Value = 111576401
$ws.Range("B15").Value = 89369
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 5447
$ws.Range("F15").Value = "Vedticka"
$ws.Range("G15").Value = "Fuscoporia viticola"
$ws.Range("H15").Value = "(Schwein.) Murrill"
$ws.Range("Q15").Value = 562964.914807545
$ws.Range("R15").Value = 6954710.791209211
$ws.Range("Z15").Value = "16:51"
$ws.Range("AB15").Value = "16:51"

# Row 16
$ws.Range("A16").Value = 111575785
$ws.Range("B16").Value = 89845
$ws.Range("E16").Value = 1209
$ws.Range("F16").Value = "Rynkskinn"
$ws.Range("G16").Value = "Phlebia centrifuga"
$ws.Range("H16").Value = "P.Karst."
$ws.Range("Q16").Value = 562859.2727272335
$ws.Range("R16").Value = 6954660.134623887
$ws.Range("Z16").Value = "16:39"
$ws.Range("AB16").Value = "16:39"

# Row 17
$ws.Range("A17").Value = 111575796
$ws.Range("B17").Value = 89686
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 658
$ws.Range("F17").Value = "Rosenticka"
$ws.Range("G17").Value = "Rhodofomes roseus"
$ws.Range("H17").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q17").Value = 562855.7640570825
$ws.Range("R17").Value = 6954651.349091855
$ws.Range("Z17").Value = "16:39"
$ws.Range("AB17").Value = "16:39"

# Row 18
$ws.Range("A18").Value = 111573569
$ws.Range("B18").Value = 96348
$ws.Range("D18").Value = "VU"
$ws.Range("E18").Value = 220787
$ws.Range("F18").Value = "Knärot"
$ws.Range("G18").Value = "Goodyera repens"
$ws.Range("H18").Value = "(L.) R. Br."
$ws.Range("Q18").Value = 562701.9737813871
$ws.Range("R18").Value = 6954788.374143652
$ws.Range("Z18").Value = "00:00"
$ws.Range("AB18").Value = "00:00"

# Row 20
$ws.Range("A20").Value = 111574240
$ws.Range("B20").Value = 56543
$ws.Range("E20").Value = 103021
$ws.Range("F20").Value = "Talltita"
$ws.Range("G20").Value = "Poecile montanus"
$ws.Range("H20").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I20").Value = "5"
$ws.Range("Q20").Value = 562533.1227179464
$ws.Range("R20").Value = 6954848.029061474
$ws.Range("Z20").Value = "15:26"
$ws.Range("AB20").Value = "15:26"
$ws.Range("AC20").Value = "Familj med 5 talltitor. Permanent revir"

# Row 21
$ws.Range("A21").Value = 111574128
$ws.Range("B21").Value = 96348
$ws.Range("D21").Value = "VU"
$ws.Range("E21").Value = 220787
$ws.Range("F21").Value = "Knärot"
$ws.Range("G21").Value = "Goodyera repens"
$ws.Range("H21").Value = "(L.) R. Br."
$ws.Range("Q21").Value = 562555.4143375416
$ws.Range("R21").Value = 6954835.60431945
$ws.Range("Z21").Value = "15:26"
$ws.Range("AB21").Value = "15:26"

# Row 22
$ws.Range("A22").Value = 111578127
$ws.Range("B22").Value = 56543
$ws.Range("E22").Value = 103021
$ws.Range("F22").Value = "Talltita"
$ws.Range("G22").Value = "Poecile montanus"
$ws.Range("H22").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("M22").Value = "lockläte, övriga läten"
$ws.Range("Q22").Value = 562937.8205991766
$ws.Range("R22").Value = 6954541.406048392
$ws.Range("Z22").Value = "18:30"
$ws.Range("AB22").Value = "18:30"

# Row 23
$ws.Range("A23").Value = 111573803
$ws.Range("Q23").Value = 562591.0245237258
$ws.Range("R23").Value = 6954847.751526525
$ws.Range("Z23").Value = "15:14"
$ws.Range("AB23").Value = "15:14"

# Row 24
$ws.Range("A24").Value = 111576771
$ws.Range("B24").Value = 96348
$ws.Range("D24").Value = "VU"
$ws.Range("E24").Value = 220787
$ws.Range("F24").Value = "Knärot"
$ws.Range("G24").Value = "Goodyera repens"
$ws.Range("H24").Value = "(L.) R. Br."
$ws.Range("M24").Value = ""
$ws.Range("Q24").Value = 562807.4867926922
$ws.Range("R24").Value = 6954821.585021482
$ws.Range("Z24").Value = "17:24"
$ws.Range("AB24").Value = "17:24"

# Row 25
$ws.Range("A25").Value = 111576450
$ws.Range("B25").Value = 96348
$ws.Range("D25").Value = "VU"
$ws.Range("E25").Value = 220787
$ws.Range("F25").Value = "Knärot"
$ws.Range("G25").Value = "Goodyera repens"
$ws.Range("H25").Value = "(L.) R. Br."
$ws.Range("Q25").Value = 562979.5212303887
$ws.Range("R25").Value = 6954739.97881452
$ws.Range("Z25").Value = "17:10"
$ws.Range("AB25").Value = "17:10"
$ws.Range("AC25").Value = "Rikligt"

# Row 26
$ws.Range("A26").Value = 111573866
$ws.Range("Q26").Value = 562601.7570288588
$ws.Range("R26").Value = 6954814.918206804
$ws.Range("Z26").Value = "15:17"
$ws.Range("AB26").Value = "15:17"

# Row 27
$ws.Range("A27").Value = 111578197
$ws.Range("Q27").Value = 563026.0554397166
$ws.Range("R27").Value = 6954541.256262898
$ws.Range("Z27").Value = "00:00"
$ws.Range("AB27").Value = "00:00"

# Row 28
$ws.Range("A28").Value = 111574689
$ws.Range("Q28").Value = 562517.0252856832
$ws.Range("R28").Value = 6954776.14289257
$ws.Range("Z28").Value = "15:47"
$ws.Range("AB28").Value = "15:47"

# Row 29
$ws.Range("A29").Value = 111575868
$ws.Range("Q29").Value = 562854.9195222461
$ws.Range("R29").Value = 6954623.341454657
$ws.Range("Z29").Value = "16:43"
$ws.Range("AB29").Value = "16:43"

# Row 30
$ws.Range("A30").Value = 111574509
$ws.Range("B30").Value = 96348
$ws.Range("E30").Value = 220787
$ws.Range("F30").Value = "Knärot"
$ws.Range("G30").Value = "Goodyera repens"
$ws.Range("H30").Value = "(L.) R. Br."
$ws.Range("Q30").Value = 562529.1073683554
$ws.Range("R30").Value = 6954769.030357216
$ws.Range("Z30").Value = "15:45"
$ws.Range("AB30").Value = "15:45"

# Row 31
$ws.Range("A31").Value = 111574403
$ws.Range("B31").Value = 89686
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 658
$ws.Range("F31").Value = "Rosenticka"
$ws.Range("G31").Value = "Rhodofomes roseus"
$ws.Range("H31").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q31").Value = 562547.0565141424
$ws.Range("R31").Value = 6954767.535469687
$ws.Range("Z31").Value = "15:26"
$ws.Range("AB31").Value = "15:26"

# Row 32
$ws.Range("A32").Value = 111576037
$ws.Range("B32").Value = 89686
$ws.Range("E32").Value = 658
$ws.Range("F32").Value = "Rosenticka"
$ws.Range("G32").Value = "Rhodofomes roseus"
$ws.Range("H32").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q32").Value = 562852.9463231879
$ws.Range("R32").Value = 6954606.325244571
$ws.Range("Z32").Value = "16:51"
$ws.Range("AB32").Value = "16:51"

# Row 33
$ws.Range("A33").Value = 111574429
$ws.Range("B33").Value = 89405
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 1202
$ws.Range("F33").Value = "Ullticka"
$ws.Range("G33").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H33").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q33").Value = 562547.0565141424
$ws.Range("R33").Value = 6954767.535469687
$ws.Range("Z33").Value = "15:42"
$ws.Range("AB33").Value = "15:42"
$ws.Range("AC33").Value = ""

